# Update countries & provincias Spain
# Applies the daily COVID data refresh: updates the "last updated" timestamp,
# refreshes several countries' case counters, and re-ranks the rows whose
# country changed position in the sorted table (Australia/Noruega,
# Guam/El Salvador, and the Republica de Africa Central / Islas Virgenes
# Britanicas / Somalia cluster).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 31 de Marzo de 2020 a las 04:50"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 164253
$ws.Range("C4").Value = 409
$ws.Range("E4").Value = 155582
$ws.Range("G4").Value = 9
$ws.Range("H4").Value = 3165

# Australia overtakes Noruega -> rows 22/23 swap labels + data
$ws.Range("A22").Value = "Australia"
$ws.Range("B22").Value = 4514
$ws.Range("C22").Value = 54
$ws.Range("D22").Value = 244
$ws.Range("E22").Value = 4251
$ws.Range("F22").Value = 28
$ws.Range("H22").Value = 19

$ws.Range("A23").Value = "Noruega"
$ws.Range("B23").Value = 4462
$ws.Range("C23").Value = 17
$ws.Range("D23").Value = 12
$ws.Range("E23").Value = 4418
$ws.Range("F23").Value = 97
$ws.Range("H23").Value = 32

# Uzbekistan (row 102)
$ws.Range("B102").Value = 150
$ws.Range("C102").Value = 1
$ws.Range("E102").Value = 141

# Trinidad y Tobago (row 118)
$ws.Range("B118").Value = 85
$ws.Range("C118").Value = 2
$ws.Range("E118").Value = 81

# El Salvador overtakes Guam -> rows 140/141 swap labels + data
$ws.Range("A140").Value = "El Salvador"
$ws.Range("C140").Value = 2
$ws.Range("E140").Value = 32
$ws.Range("F140").Value = 5
$ws.Range("H140").Value = 0

$ws.Range("A141").Value = "Guam"
$ws.Range("B141").Value = 32
$ws.Range("E141").Value = 31
$ws.Range("F141").Value = 0
$ws.Range("H141").Value = 1

# Republica de Africa Central / Islas Virgenes Britanicas / Somalia re-rank
$ws.Range("A198").Value = "Islas Virgenes Britanicas"
$ws.Range("C198").Value = 1

$ws.Range("A199").Value = "Republica de Africa Central"

$ws.Range("A202").Value = "Somalia"
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 1
$ws.Range("E202").Value = 2
